# Regenerate save_data to use K instead of Strike# for column G (rows 2-17).
# This updates the "K" column values on the active sheet to their recalculated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 1
    10 = 3
    11 = 0
    12 = 0
    14 = 1
    16 = 0
    17 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
